# Testdata correction due to LIVEHTA-2596 and LIVEHTA-2980
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- scenario1 block (rows 2-8): corrected ExpectedPrismaCount values ---
$ws.Range("J2").Value = 14
$ws.Range("J3").Value = 4
$ws.Range("J4").Value = 4
$ws.Range("J5").Value = 2
$ws.Range("J7").Value = 5
$ws.Range("J8").Value = 4

# --- scenario2 block (rows 10-16): corrected ExpectedPrismaCount values ---
$ws.Range("J10").Value = 14
$ws.Range("J11").Value = 6

# New section row for scenario2 - an extra study_design section/checkbox pairing
# that was missing from row 13, plus its corrected ExpectedPrismaCount.
$ws.Range("F13").Value = "study_design_section2_1"
$ws.Range("G13").Value = "study_design_section2_1_checkbox"
$ws.Range("H13").Value = "study_design_section"
$ws.Range("J13").Value = 2

$ws.Range("J14").Value = 5
$ws.Range("J16").Value = 5

# --- Restore the view/selection state (scroll position + active cell) ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 3
$win.ScrollRow = 1
$ws.Range("C17").Select()
